$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on column A (Beteckning)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C ("Förändrad") holds a date serial value that changed from
# 45186 (2023-09-17) to 45188 (2023-09-19) for every data row (2..lastRow).
$ws.Range("C2:C$lastRow").Value = 45188
